# Rename the three header cells (B1, C1, D1) on Sheet1 to match the new
# naming convention used by the services py/bat setup.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "unique_id"
$ws.Range("C1").Value = "github_url"
$ws.Range("D1").Value = "highlight_app_id"

# Widen column D so the longer "highlight_app_id" header fits.
$ws.Columns.Item(4).ColumnWidth = 15.43

# Reflect the author's active selection at save time.
$ws.Range("D13").Select()
